# Apply the changes described by the commit:
# - Update the "Date" metadata value
# - Swap the "Mapping: RIM Mapping" / "Mapping: Spécification métier vers
#   l'extension ROR AvailableTimeTypeOfTime" columns (header + data) on the
#   Elements sheet, so the French mapping column now comes before the RIM
#   mapping column
# - Best-effort match of the resulting (auto best-fit) column widths

$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" value ---------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Elements sheet: swap the AK / AL (Mapping) columns --------------
$els = $wb.Worksheets.Item("Elements")

# Find the last used row on the Elements sheet
$lastRow = $els.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $els.Cells.Item($r, 37)   # column AK
    $alCell = $els.Cells.Item($r, 38)   # column AL

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# --- 3. Best-effort column width swap (bestFit auto-resize) -------------
# The workbook columns are "best fit" so their width simply follows the
# (now-swapped) column content; reproduce that as closely as the
# ColumnWidth setter's precision allows.
$els.Columns.Item(37).ColumnWidth = 81.15
$els.Columns.Item(38).ColumnWidth = 24.15
